# BOT; UPDATE DATA
# Adds the 2020-04-19 (serial 43940) daily/cumulative consultation counts as a
# new row 85, and pushes the existing footnote row down to row 86.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the footnote text currently stored in B85 so it can be moved to B86.
$footnote = $ws.Range("B85").Value2

# --- Row 85 becomes the new data row for 2020-04-19 ---
$ws.Range("A85").Value = 43940
$ws.Range("B85").Value = 523
$ws.Range("C85").Value = 25980
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 5665

# --- Row 86 becomes the (shifted down) footnote row ---
$ws.Range("B86").Value = $footnote

# Update the active selection to match the edited cell (E85), as Excel would
# leave the cursor there after entering the last value of the new row.
$ws.Range("E85").Select() | Out-Null
